$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046410094356494
$ws.Range("D2").Value = 1.052502357523927
$ws.Range("E2").Value = 1.0535431997867
$ws.Range("F2").Value = 1.063154658258427
$ws.Range("I2").Value = 1.046370944666637
$ws.Range("J2").Value = 1.051464639180589
$ws.Range("K2").Value = 1.055250912285313
$ws.Range("L2").Value = 1.056288878300932
$ws.Range("M2").Value = 1.065874064711381

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047272777979757
$ws.Range("D3").Value = 1.053174587009434
$ws.Range("E3").Value = 1.054300472242774
$ws.Range("F3").Value = 1.063969289289568
$ws.Range("I3").Value = 1.046599402790632
$ws.Range("J3").Value = 1.051975807316008
$ws.Range("K3").Value = 1.055736361232915
$ws.Range("L3").Value = 1.056859359196669
$ws.Range("M3").Value = 1.066503652169136

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.047831598638368
$ws.Range("D4").Value = 1.053610051471106
$ws.Range("E4").Value = 1.054791384495719
$ws.Range("F4").Value = 1.064497371788746
$ws.Range("I4").Value = 1.046746182611589
$ws.Range("J4").Value = 1.052306500533726
$ws.Range("K4").Value = 1.056050282978239
$ws.Range("L4").Value = 1.057228738463213
$ws.Range("M4").Value = 1.066911343259094

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048066670468117
$ws.Range("D5").Value = 1.05379323556001
$ws.Range("E5").Value = 1.054997979158089
$ws.Range("F5").Value = 1.064719605997859
$ws.Range("I5").Value = 1.046807637173305
$ws.Range("J5").Value = 1.052445506687237
$ws.Range("K5").Value = 1.0561822073908
$ws.Range("L5").Value = 1.057384081523131
$ws.Range("M5").Value = 1.067082808249699

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048106148438105
$ws.Range("D6").Value = 1.053823999656794
$ws.Range("E6").Value = 1.055032679874371
$ws.Range("F6").Value = 1.064756933432829
$ws.Range("I6").Value = 1.046817940882627
$ws.Range("J6").Value = 1.052468845381367
$ws.Range("K6").Value = 1.05620435521002
$ws.Range("L6").Value = 1.057410167549512
$ws.Range("M6").Value = 1.06711160213218

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04783473911733
$ws.Range("D7").Value = 1.05361249873689
$ws.Range("E7").Value = 1.054794144180609
$ws.Range("F7").Value = 1.064500340398368
$ws.Range("I7").Value = 1.046747004760378
$ws.Range("J7").Value = 1.052308358010257
$ws.Range("K7").Value = 1.056052045950098
$ws.Range("L7").Value = 1.057230813946479
$ws.Range("M7").Value = 1.06691363410311

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.046701515740858
$ws.Range("D8").Value = 1.052729438720105
$ws.Range("E8").Value = 1.0537989347974
$ws.Range("F8").Value = 1.063429766339915
$ws.Range("I8").Value = 1.046448369635741
$ws.Range("J8").Value = 1.051637403886412
$ws.Range("K8").Value = 1.055415011856546
$ws.Range("L8").Value = 1.056481624438279
$ws.Range("M8").Value = 1.066086772452097

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.044709349095275
$ws.Range("D9").Value = 1.051177182344819
$ws.Range("E9").Value = 1.052052270941793
$ws.Range("F9").Value = 1.061550728715294
$ws.Range("I9").Value = 1.045914147078352
$ws.Range("J9").Value = 1.050454643086296
$ws.Range("K9").Value = 1.054291032292907
$ws.Range("L9").Value = 1.055163358215764
$ws.Range("M9").Value = 1.064632149204381

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043384505847398
$ws.Range("D10").Value = 1.050145010511039
$ws.Range("E10").Value = 1.050892658676792
$ws.Range("F10").Value = 1.060303160720003
$ws.Range("I10").Value = 1.045552675682619
$ws.Range("J10").Value = 1.049665907856401
$ws.Range("K10").Value = 1.053540820388045
$ws.Range("L10").Value = 1.054285877984772
$ws.Range("M10").Value = 1.063664116507995

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.042811629307164
$ws.Range("D11").Value = 1.049698723224958
$ws.Range("E11").Value = 1.050391702453033
$ws.Range("F11").Value = 1.059764189429707
$ws.Range("I11").Value = 1.045394902222928
$ws.Range("J11").Value = 1.049324338948396
$ws.Range("K11").Value = 1.053215775317032
$ws.Range("L11").Value = 1.053906260887236
$ws.Range("M11").Value = 1.063245374531637

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.042598957648596
$ws.Range("D12").Value = 1.049533051697716
$ws.Range("E12").Value = 1.050205801679255
$ws.Range("F12").Value = 1.059564178826633
$ws.Range("I12").Value = 1.045336110543973
$ws.Range("J12").Value = 1.049197460136542
$ws.Range("K12").Value = 1.053095010631725
$ws.Range("L12").Value = 1.053765306371852
$ws.Range("M12").Value = 1.063089900200524

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.042644570976323
$ws.Range("D13").Value = 1.049568584273913
$ws.Range("E13").Value = 1.050245670001556
$ws.Range("F13").Value = 1.059607073271414
$ws.Range("I13").Value = 1.045348730034681
$ws.Range("J13").Value = 1.049224676273729
$ws.Range("K13").Value = 1.053120916333811
$ws.Range("L13").Value = 1.05379553920946
$ws.Range("M13").Value = 1.063123247001521

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.042794047340253
$ws.Range("D14").Value = 1.049685026725867
$ws.Range("E14").Value = 1.050376332222235
$ws.Range("F14").Value = 1.05974765265805
$ws.Range("I14").Value = 1.045390046306562
$ws.Range("J14").Value = 1.049313851198888
$ws.Range("K14").Value = 1.053205793443634
$ws.Range("L14").Value = 1.053894608470234
$ws.Range("M14").Value = 1.063232521638011

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.042886160636684
$ws.Range("D15").Value = 1.04975678399432
$ws.Range("E15").Value = 1.050456860995384
$ws.Range("F15").Value = 1.059834293139596
$ws.Range("I15").Value = 1.045415477791653
$ws.Range("J15").Value = 1.049368794201581
$ws.Range("K15").Value = 1.053258085299052
$ws.Range("L15").Value = 1.053955655268876
$ws.Range("M15").Value = 1.063299858020017

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043422542544087
$ws.Range("D16").Value = 1.050174642975127
$ws.Range("E16").Value = 1.050925930165103
$ws.Range("F16").Value = 1.060338956632884
$ws.Range("I16").Value = 1.04556312023646
$ws.Range("J16").Value = 1.049688575899783
$ws.Range("K16").Value = 1.053562388493001
$ws.Range("L16").Value = 1.054311079158752
$ws.Range("M16").Value = 1.063691916073266

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.043759213203461
$ws.Range("D17").Value = 1.050436930148016
$ws.Range("E17").Value = 1.051220477691646
$ws.Range("F17").Value = 1.060655850435328
$ws.Range("I17").Value = 1.045655397286463
$ws.Range("J17").Value = 1.049889156358778
$ws.Range("K17").Value = 1.053753217677048
$ws.Range("L17").Value = 1.054534118593506
$ws.Range("M17").Value = 1.063937957993481

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.043955663553817
$ws.Range("D18").Value = 1.050589980406753
$ws.Range("E18").Value = 1.051392394420273
$ws.Range("F18").Value = 1.060840808149807
$ws.Range("I18").Value = 1.045709099885684
$ws.Range("J18").Value = 1.050006147470151
$ws.Range("K18").Value = 1.053864505883019
$ws.Range("L18").Value = 1.05466424620164
$ws.Range("M18").Value = 1.064081510747569

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.044022660897209
$ws.Range("D19").Value = 1.05064217715117
$ws.Range("E19").Value = 1.051451032546115
$ws.Range("F19").Value = 1.060903894062277
$ws.Range("I19").Value = 1.045727390528874
$ws.Range("J19").Value = 1.050046037684732
$ws.Range("K19").Value = 1.053902448983391
$ws.Range("L19").Value = 1.054708621818111
$ws.Range("M19").Value = 1.064130465375822

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.043723083750033
$ws.Range("D20").Value = 1.050408782729849
$ws.Range("E20").Value = 1.051188863917025
$ws.Range("F20").Value = 1.060621838414723
$ws.Range("I20").Value = 1.04564550934654
$ws.Range("J20").Value = 1.049867636386094
$ws.Range("K20").Value = 1.053732745490892
$ws.Range("L20").Value = 1.05451018520301
$ws.Range("M20").Value = 1.063911555806288

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.042750026945753
$ws.Range("D21").Value = 1.049650734566517
$ws.Range("E21").Value = 1.050337850562918
$ws.Range("F21").Value = 1.059706250345568
$ws.Range("I21").Value = 1.045377884866258
$ws.Range("J21").Value = 1.049287591531174
$ws.Range("K21").Value = 1.053180800025628
$ws.Range("L21").Value = 1.053865433590174
$ws.Range("M21").Value = 1.06320034117172

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.04213892372673
$ws.Range("D22").Value = 1.049174695826753
$ws.Range("E22").Value = 1.049803807786457
$ws.Range("F22").Value = 1.059131668597817
$ws.Range("I22").Value = 1.045208533787225
$ws.Range("J22").Value = 1.048922866539275
$ws.Range("K22").Value = 1.052833605924925
$ws.Range("L22").Value = 1.053460355259775
$ws.Range("M22").Value = 1.062753549645546

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04246281457466
$ws.Range("D23").Value = 1.049426997750402
$ws.Range("E23").Value = 1.050086816366983
$ws.Range("F23").Value = 1.059436161773842
$ws.Range("I23").Value = 1.045298412570522
$ws.Range("J23").Value = 1.04911621635604
$ws.Range("K23").Value = 1.053017675196038
$ws.Range("L23").Value = 1.053675065779639
$ws.Range("M23").Value = 1.062990365986272

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.043739408860064
$ws.Range("D24").Value = 1.050421501143367
$ws.Range("E24").Value = 1.051203148474722
$ws.Range("F24").Value = 1.060637206615606
$ws.Range("I24").Value = 1.045649977655065
$ws.Range("J24").Value = 1.049877360348187
$ws.Range("K24").Value = 1.053741996050175
$ws.Range("L24").Value = 1.054520999571798
$ws.Range("M24").Value = 1.063923485692762

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045223802978284
$ws.Range("D25").Value = 1.051578015673089
$ws.Range("E25").Value = 1.052502982012899
$ws.Range("F25").Value = 1.06203561025358
$ws.Range("I25").Value = 1.046053198193782
$ws.Range("J25").Value = 1.050760460945028
$ws.Range("K25").Value = 1.054581770971372
$ws.Range("L25").Value = 1.055503927492445
$ws.Range("M25").Value = 1.065007908582331
